# Update workbook strings to reflect the new build/version info.
# Old version string: "mines - January 30 (built on February 02 2026 12.49.33 EST)"
# New version string: "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

$wb = $excel.ActiveWorkbook

$newVersion = "mines - version 1.0.0 (Feb 3 2026) (built on February 03 2026 10.14.00 EST)"

# --- "About" worksheet ---
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("A2").Value = "Version: $newVersion"

$wsAbout.Range("A6").Value = "Recommended Citation:  `"Global Energy Monitor, Coal mine boundaries and methane sources for Shaqu No.2 Coal Mine, China, M1195, version '$newVersion'. (See the CC license for attribution requirements if sharing or adapting the data set.)"

# --- "Boundaries and methane sources" worksheet ---
$wsData = $wb.Worksheets.Item("Boundaries and methane sources")

# Column S holds "build_version" (header in row 1); data rows 2-10 need updating.
for ($r = 2; $r -le 10; $r++) {
    $wsData.Cells.Item($r, 19).Value = $newVersion
}
